$wb = $excel.ActiveWorkbook

# Previously active sheet was Slovakia - switch to the Turkey sheet and
# select all its cells (as if via the corner "select all" button) before
# duplicating it.
$turkey = $wb.Worksheets.Item("Turkey")
$turkey.Activate()
$turkey.Cells.Select() | Out-Null

# Duplicate the Turkey sheet as a template for the new Croatia market sheet,
# placing it after the last sheet.
$turkey.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$croatia = $wb.Worksheets.Item($wb.Worksheets.Count)
$croatia.Name = "Croatia"

# Fill in market-specific values.
$croatia.Range("B2").Value = "Croatia Market"
$croatia.Range("B4").Value = "NGC-3139/T2415"

# Leave the new sheet active with B4 selected, matching the edit session end state.
$croatia.Activate()
$croatia.Range("B4").Select() | Out-Null
